$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

$ws.Range("A5").Value = ""
$ws.Range("B5").Value = "mc"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "Wat betekend gelijkstroom3"
$ws.Range("E5").Value = "['test', 'test1', 'test2']"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
